$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = $origStyle
}

# Row 2 - date only
Set-TextValue "A2" "2025/11/28"

# Row 8 - date + EBITDA
Set-TextValue "A8" "2025/11/28"
Set-TextValue "B8" "7.61"

# Row 14 - date + EBITDA
Set-TextValue "A14" "2025/11/28"
Set-TextValue "B14" "2.82"

# Row 20 - date only
Set-TextValue "A20" "2025/11/28"

# Row 26 - date only
Set-TextValue "A26" "2025/11/28"

# Row 32 - date only
Set-TextValue "A32" "2025/11/28"

# Row 38 - date only
Set-TextValue "A38" "2025/11/28"

# Row 44 - date only
Set-TextValue "A44" "2025/11/28"

# Row 50 - date only
Set-TextValue "A50" "2025/11/28"

# Row 56 - date + EBITDA
Set-TextValue "A56" "2025/11/28"
Set-TextValue "B56" "35.17"

# Row 62 - date only
Set-TextValue "A62" "2025/11/28"

# Row 68 - date only
Set-TextValue "A68" "2025/11/28"

# Row 74 - date only
Set-TextValue "A74" "2025/11/28"
